$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The journal table body (row 3) is the template row; clone its formatting
# (style "9" on A:G) down through the newly-appended rows 4-19.
$ws.Range("A3:G3").Copy()
$ws.Range("A4:G19").PasteSpecial(-4122)

$cols = @("B", "C", "D", "E", "F", "G")
$oddValues = @("Тимків Віталій Дмитрович", "А0000", "Тимків Дмитро Віталійович", "Уганда, гасити вагнерів", "01.01.2025 Краківець", "01.01.2026 Подобовець")
$evenValues = @("Тимків Віталій Дмитрович 1", "А0000 1", "Тимків Дмитро Віталійович 1", "Уганда, гасити вагнерів 1", "01.01.2025 Краківець 1", "01.01.2026 Подобовець 1")

for ($i = 1; $i -le 16; $i++) {
    $row = $i + 3
    $ws.Cells.Item($row, 1).Value = $i
    if ($i % 2 -eq 1) {
        $rowValues = $oddValues
    } else {
        $rowValues = $evenValues
    }
    for ($j = 0; $j -lt 6; $j++) {
        $ws.Cells.Item($row, $cols[$j]).Value = $rowValues[$j]
    }
}
